$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vokabeltrainer")
$ws.Activate()

$ws.Range("E4").Value = 13
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 1
$ws.Range("E11").Value = 8
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 5
$ws.Range("E15").Value = 3
$ws.Range("E16").Value = 2
$ws.Range("E17").Value = 8
$ws.Range("E18").Value = 5
$ws.Range("E19").Value = 5
$ws.Range("E20").Value = 3

$ws.Range("F6").Select()
